$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9999958875889788
$ws.Range("C2").Value = 0.2893471185498306
$ws.Range("D2").Value = 0.7106487690391482
$ws.Range("E2").Value = 0.06651617804570137
$ws.Range("F2").Value = 27.65080570750881
$ws.Range("H2").Value = 2.41

$ws.Range("B3").Value = 0.9999919561955879
$ws.Range("C3").Value = 0.3442991438669198
$ws.Range("D3").Value = 0.6556928123286681
$ws.Range("E3").Value = 0.09302716390936586
$ws.Range("F3").Value = 26.56023451251171
$ws.Range("H3").Value = 1.13

$ws.Range("B4").Value = 0.9999786254387332
$ws.Range("C4").Value = 0.3155303178371341
$ws.Range("D4").Value = 0.684448307601599
$ws.Range("E4").Value = 0.1516448390638694
$ws.Range("F4").Value = 27.13664405141674
$ws.Range("H4").Value = 2.93

$ws.Range("C5").Value = 0.5550794535011663
$ws.Range("D5").Value = 0.4016976419844882
$ws.Range("E5").Value = 6.819240442267574
$ws.Range("F5").Value = 21.87863690959999
$ws.Range("H5").Value = 2.34

$ws.Range("H6").Value = 1.67

$ws.Range("C7").Value = -1.198298859569846
$ws.Range("D7").Value = 1.369583775334443
$ws.Range("F7").Value = 48.63202988343252
$ws.Range("H7").Value = 2.47

$ws.Range("C8").Value = 0.2277983465155334
$ws.Range("D8").Value = 0.7703424069921297
$ws.Range("E8").Value = 1.414319454662004
$ws.Range("F8").Value = 28.82334589671151
